# Apply the cryptos-list refresh described by the commit diff.
# Price cells whose text would otherwise be auto-parsed as a number by Excel
# (single-dot decimals, e.g. "331.50") are written with a leading apostrophe,
# matching Excel's own text-prefix convention, so they keep round-tripping as
# plain text instead of being coerced to a numeric value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.463.62"
$ws.Range("E2").Value = "  -0.44%  "
$ws.Range("D3").Value = "1.824.05"
$ws.Range("E3").Value = "  -1.92%  "
$ws.Range("E4").Value = "  -0.50%  "
$ws.Range("D5").Value = "'331.50"
$ws.Range("E5").Value = "  -0.74%  "
$ws.Range("E6").Value = "  -0.58%  "
$ws.Range("D7").Value = "'0.4590"
$ws.Range("E7").Value = "  -1.86%  "
$ws.Range("D8").Value = "'0.3809"
$ws.Range("E8").Value = "  -2.13%  "
$ws.Range("D9").Value = "'46.47"
$ws.Range("E9").Value = "  +2.48%  "
$ws.Range("D10").Value = "'0.07895"
$ws.Range("E10").Value = "  -0.99%  "
$ws.Range("D11").Value = "'0.9679"
$ws.Range("E11").Value = "  -3.15%  "
$ws.Range("E12").Value = "  -2.84%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'5.884"
$ws.Range("E13").Value = "  -1.41%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.829.39"
$ws.Range("E14").Value = "  -1.90%  "
$ws.Range("D15").Value = "'7.047"
$ws.Range("E15").Value = "  -2.66%  "
$ws.Range("E16").Value = "  -0.75%  "
$ws.Range("D17").Value = "'89.80"
$ws.Range("E17").Value = "  +2.07%  "
$ws.Range("D18").Value = "'0.06617"
$ws.Range("E18").Value = "  -1.62%  "
$ws.Range("E19").Value = "  -1.51%  "
$ws.Range("E20").Value = "  +0.79%  "
$ws.Range("D21").Value = "'1.005"
$ws.Range("E21").Value = "  -0.50%  "
$ws.Range("D22").Value = "27.441.37"
$ws.Range("E22").Value = "  -0.48%  "
$ws.Range("D23").Value = "'5.332"
$ws.Range("E23").Value = "  -1.81%  "
$ws.Range("D24").Value = "'10.80"
$ws.Range("E24").Value = "  -0.33%  "
$ws.Range("D25").Value = "'2.303"
$ws.Range("E25").Value = "  -0.13%  "
$ws.Range("D26").Value = "2.059.99"
$ws.Range("E26").Value = "  -1.17%  "
$ws.Range("D27").Value = "'155.72"
$ws.Range("E27").Value = "  -1.92%  "
$ws.Range("D28").Value = "'19.35"
$ws.Range("E28").Value = "  -1.86%  "
$ws.Range("D29").Value = "'2.056"
$ws.Range("E29").Value = "  -3.15%  "
$ws.Range("D30").Value = "'5.270"
$ws.Range("E30").Value = "  -2.18%  "
$ws.Range("D31").Value = "'118.24"
$ws.Range("E31").Value = "  -2.40%  "
$ws.Range("D32").Value = "'0.9443"
$ws.Range("E32").Value = "  -2.74%  "
$ws.Range("D33").Value = "'0.09325"
$ws.Range("E33").Value = "  -1.31%  "
$ws.Range("D34").Value = "'3.595"
$ws.Range("E34").Value = "  -1.32%  "
$ws.Range("D35").Value = "'5.238"
$ws.Range("E35").Value = "  -0.73%  "
$ws.Range("D36").Value = "'1.319"
$ws.Range("E36").Value = "  -0.45%  "
$ws.Range("D37").Value = "'0.05922"
$ws.Range("E37").Value = "  -1.81%  "
$ws.Range("D38").Value = "'0.02183"
$ws.Range("E38").Value = "  -1.63%  "
$ws.Range("D39").Value = "'1.160"
$ws.Range("E39").Value = "  -2.48%  "
$ws.Range("D40").Value = "'8.026"
$ws.Range("E40").Value = "  -2.33%  "
$ws.Range("D41").Value = "'0.5750"
$ws.Range("E41").Value = "  -2.65%  "
$ws.Range("E42").Value = "  -2.44%  "
$ws.Range("D43").Value = "'10.02"
$ws.Range("E43").Value = "  -1.89%  "
$ws.Range("D44").Value = "'1.261"
$ws.Range("E44").Value = "  +0.81%  "
$ws.Range("D45").Value = "'11.94"
$ws.Range("E45").Value = "  -1.61%  "
$ws.Range("E46").Value = "  -3.17%  "
$ws.Range("D47").Value = "'1.866"
$ws.Range("E47").Value = "  -2.49%  "
$ws.Range("D48").Value = "'0.06607"
$ws.Range("E48").Value = "  -2.27%  "
$ws.Range("D49").Value = "'110.44"
$ws.Range("E49").Value = "  -1.74%  "
$ws.Range("E50").Value = "  -0.72%  "
$ws.Range("D51").Value = "'1.041"
$ws.Range("E51").Value = "  -1.23%  "
